$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Car"

# Reposition the workbook window (best-effort; harmless if unsupported)
try {
    $win = $excel.ActiveWindow
    $win.Left = 1020
    $win.Top = 1065
} catch {
}

# --- Column A: single letter labels (rows 2-8) ---
$ws.Range("A2").Value = "a"
$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "c"
$ws.Range("A5").Value = "d"
$ws.Range("A6").Value = "e"
$ws.Range("A7").Value = "f"
$ws.Range("A8").Value = "g"

# --- Row 2 ---
$ws.Range("B2").Value = "b,3,0,0"
$ws.Range("C2").Value = "c,7,0,0"

# --- Row 3 ---
$ws.Range("B3").Value = "a,3,0,0"
$ws.Range("C3").Value = "d,5,0,0"

# --- Row 4 ---
$ws.Range("B4").Value = "a,7,0,0"
$ws.Range("C4").Value = "d,9,0,0"

# --- Row 5 ---
$ws.Range("B5").Value = "b,5,0,0"
$ws.Range("C5").Value = "c,9,0,0"
$ws.Range("D5").Value = "e,2,0,0"
$ws.Range("E5").Value = "f,1,0,0"

# --- Row 6 ---
$ws.Range("B6").Value = "d,2,0,0"
$ws.Range("C6").Value = "g,3,0,0"

# --- Row 7 ---
$ws.Range("B7").Value = "d,1,0,0"
$ws.Range("C7").Value = "g,4,0,0"

# --- Row 8 ---
$ws.Range("B8").Value = "e,3,0,0"
$ws.Range("C8").Value = "f,4,0,0"

# Center-align the whole populated block, including the blank D/E cells
# on rows that don't otherwise have values (2,3,4,6,7,8)
$ws.Range("A2:E8").HorizontalAlignment = -4108

# Restore selection / active cell
$ws.Range("G13").Select()
